$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 24-28 use the "General + border" style that currently only
# appears on A1 (style index 2 in the original cellXfs table). Copy that
# cell's format onto the whole new block first, then fill in the values -
# this reproduces the s="2" styling on every new cell exactly.
$fmtSrc = $ws.Range("A1")
$newBlock = $ws.Range("A24:F28")
$fmtSrc.Copy()
$newBlock.PasteSpecial(-4122)  # xlPasteFormats

# Match the row heights Excel stamps on newly-typed rows (15.75, explicit).
$ws.Range("A24:F28").RowHeight = 15.75

# Row 24: one new [E, A] ballot
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "E"
$ws.Range("C24").Value = "A"

# Rows 25-28: four new [A] ballots
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "A"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "A"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "A"

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "A"
